$d = $word.ActiveDocument

# Split the centered name line into two paragraphs: the existing
# "Dheeraj Chand" heading line, followed by a brand-new centered
# contact-info line. Using Find/Replace with the paragraph-mark wildcard
# (^p) lets Word insert a genuine new paragraph that inherits the
# original paragraph's formatting (centered) while keeping the new run
# free of the name run's bold/size-28 character formatting.
$found = $d.Content.Find.Execute(
    "Dheeraj Chand", $true, $false, $false, $false, $false, $true, 1, $false,
    "Dheeraj Chand^p202.550.7110 | dheeraj.chand@gmail.com | https://www.dheerajchand.com | https://www.linkedin.com/in/dheerajchand/ | Austin, TX",
    2
)
